# Update the "Förändrad" (Changed) date column (C) from 45170 to 45174
# for all data rows (2 through 99) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 99 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45170) {
        $cell.Value = 45174
    }
}
